# input gain calculation for ADC
#
# 1) Rename the existing sheet "Sheet1" -> "RCL" and tweak a couple of its
#    input cells (B1, B4). Formulas on that sheet recalc automatically.
# 2) Add a brand-new sheet "Sheet2" (placed right after "RCL", and becomes
#    the active tab) that performs the ADC input-gain calculation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "RCL" (was "Sheet1")
# ---------------------------------------------------------------------
$rcl = $wb.Worksheets.Item(1)
$rcl.Name = "RCL"

$rcl.Range("B1").Value = 10
$rcl.Range("B4").Value = 2

$rcl.Columns.Item(1).ColumnWidth = 14.140625
$rcl.Columns.Item(5).ColumnWidth = 12.85546875

$null = $rcl.Range("B1").Select()
$excel.ActiveWindow.Zoom = 145

# ---------------------------------------------------------------------
# New sheet "Sheet2"
# ---------------------------------------------------------------------
$s2 = $wb.Worksheets.Add($null, $rcl)
$s2.Name = "Sheet2"

$s2.Range("A1").Value = "Vout"
$s2.Range("B1").Value = 3

$s2.Range("A2").Value = "Vin"
$s2.Range("B2").Value = 2.5

$s2.Range("A3").Value = "Gain"
$s2.Range("B3").Formula = '=B1/B2'

# NOTE: "Rf" (row 6) is entered before "Rg" (row 4) so the shared-string
# table allocates indices in the same order as the authored workbook.
$s2.Range("A6").Value = "Rf"
$s2.Range("A4").Value = "Rg"
$s2.Range("B4").Value = 330

$s2.Range("B6").Formula = '=B1 * B4 / B2'
$s2.Range("D6").Value = 392
$s2.Range("E6").Value = 402
$s2.Range("F6").Formula = '=768 / 2'

$s2.Range("D7").Formula = '=D6/$B$4'
$s2.Range("E7").Formula = '=E6/$B$4'
$s2.Range("F7").Formula = '=F6/$B$4'

$s2.Range("D8").Formula = '=$B$2 * D7'
$s2.Range("E8").Formula = '=$B$2 * E7'
$s2.Range("F8").Formula = '=$B$2 * F7'

$s2.Range("A11").Value = "Rin"
$s2.Range("B11").Formula = '=B4/(1-F6/(2*(B4+F6)))'

$s2.Range("A12").Value = "Imp"
$s2.Range("B12").Value = 50

$s2.Range("A13").Value = "Rt"
$s2.Range("B13").Formula = '=1/ABS(1/B12-1/B11)'
$s2.Range("D13").Value = 56.2

$s2.Range("D15").Formula = '= 1 / (1 / B11 + 1 / D13)'

# View state: Sheet2 becomes the active/visible tab, zoomed in, scrolled
# down a bit with the selection parked on E14.
$null = $s2.Activate()
$excel.ActiveWindow.Zoom = 250
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$null = $s2.Range("E14").Select()
